# RN01.docx: add functional requirements RN02-RN07 to the requirements table,
# change the RN01 description, and grow the table to 3 columns (new "depends
# on" column on the right).
#
# Word's Range.InsertXML() happens to mis-nest a replacement <w:tbl> when the
# target range is scoped to the table itself (it ends up inserting the new
# table inside the first cell of the old one). Selecting a range that spans
# the whole document body - table plus the trailing empty paragraph that
# already follows it - avoids that and performs a clean top-level
# replacement, so we rebuild the table (plus that trailing paragraph,
# unchanged) as one XML blob and drop it in through $d.Content.

$w = "http://schemas.openxmlformats.org/wordprocessingml/2006/main"

function Tc([string]$widthDxa, [string]$innerXml) {
    $open = "<w:tc><w:tcPr><w:tcW w:w=`"$widthDxa`" w:type=`"dxa`"/></w:tcPr><w:p>"
    $close = "</w:p></w:tc>"
    $result = $open + $innerXml + $close
    return $result
}

function SimpleRun([string]$text) {
    $result = "<w:r><w:t>$text</w:t></w:r>"
    return $result
}

$col1Width = "879"
$col2Width = "4928"
$col3Width = "992"

$ids = @("[RN01]", "[RN02]", "[RN03]", "[RN04]", "[RN05]", "[RN06]", "[RN07]")

$col2Texts = @(
    "O sistema deve permitir que o usuário crie seu cadastro.",
    "",
    "O sistema deve permitir que o usuário altere seu cadastro.",
    "O sistema deve permitir que o usuário crie uma nova despesa.",
    "O sistema deve exigir que o usuário insira os participantes da despesa.",
    "O sistema deve permitir que o usuário insira os itens da despesa.",
    "O sistema deve permitir que o usuário veja em seu perfil o saldo com cada morador da casa."
)

# RN02's paragraph has the word "login" split into its own run and flagged
# with proofErr spell-check markers, so it is built separately.
$rn02Inner = "<w:r><w:t xml:space=`"preserve`">O sistema deve permitir que o usuário faça o </w:t></w:r>" +
    "<w:proofErr w:type=`"spellStart`"/><w:r><w:t>login</w:t></w:r><w:proofErr w:type=`"spellEnd`"/>" +
    "<w:r><w:t xml:space=`"preserve`"> com os dados cadastrados.</w:t></w:r>"

$col3Texts = @("", "", "", "", "[RN04]", "[RN04]", "")

$rowsXml = ""
for ($i = 0; $i -lt $ids.Count; $i++) {
    $id = $ids[$i]
    $col2Text = $col2Texts[$i]
    $col3Text = $col3Texts[$i]

    $idRun = SimpleRun $id
    $col1Cell = Tc $col1Width $idRun

    if ($i -eq 1) {
        $col2Inner = $rn02Inner
    } else {
        $col2Inner = SimpleRun $col2Text
    }
    $col2Cell = Tc $col2Width $col2Inner

    $col3Inner = ""
    if ($col3Text -ne "") {
        $col3Inner = SimpleRun $col3Text
    }
    $col3Cell = Tc $col3Width $col3Inner

    $row = "<w:tr>" + $col1Cell + $col2Cell + $col3Cell + "</w:tr>"
    $rowsXml = $rowsXml + $row
}

# Final row: empty except for the _GoBack bookmark preserved from the
# original last paragraph of the table.
$bookmarkInner = '<w:bookmarkStart w:id="0" w:name="_GoBack"/><w:bookmarkEnd w:id="0"/>'
$lastCol1 = Tc $col1Width $bookmarkInner
$lastCol2 = Tc $col2Width ""
$lastCol3 = Tc $col3Width ""
$lastRow = "<w:tr>" + $lastCol1 + $lastCol2 + $lastCol3 + "</w:tr>"
$rowsXml = $rowsXml + $lastRow

$tblPr = "<w:tblPr><w:tblStyle w:val=`"Tabelacomgrade`"/><w:tblW w:w=`"6799`" w:type=`"dxa`"/>" +
    "<w:tblLook w:val=`"04A0`" w:firstRow=`"1`" w:lastRow=`"0`" w:firstColumn=`"1`" w:lastColumn=`"0`" w:noHBand=`"0`" w:noVBand=`"1`"/></w:tblPr>"
$tblGrid = "<w:tblGrid><w:gridCol w:w=`"$col1Width`"/><w:gridCol w:w=`"$col2Width`"/><w:gridCol w:w=`"$col3Width`"/></w:tblGrid>"
$tblOpen = "<w:tbl xmlns:w=`"$w`">" + $tblPr + $tblGrid
$tblClose = "</w:tbl>"
$tblXml = $tblOpen + $rowsXml + $tblClose

# Preserve the paragraph that already followed the table, unchanged.
$trailingParagraph = '<w:p w:rsidR="00C95E3C" w:rsidRDefault="00C95E3C"/>'

$newXml = $tblXml + $trailingParagraph

$d = $word.ActiveDocument
$d.Content.InsertXML($newXml)
